$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.912.03"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.359.94"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.672"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.05"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.27"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  +2.67%  "

$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.78"
$ws.Range("E11").Value = "  +6.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "36.33"
$ws.Range("E12").Value = "  +12.76%  "

$ws.Range("E13").Value = "  +0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.28"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.34"
$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("D17").Value = "2.365.39"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "43.814.18"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("E19").Value = "  +1.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.98"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "253.87"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("E24").Value = "  +3.35%  "

$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.51"
$ws.Range("E27").Value = "  -2.13%  "

$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.34"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  +0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.134"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0751"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.08"
$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.38"
$ws.Range("E35").Value = "  -0.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.81"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.61"
$ws.Range("E37").Value = "  +5.13%  "

$ws.Range("E38").Value = "  +2.86%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.45"
$ws.Range("E40").Value = "  +13.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.33"
$ws.Range("E41").Value = "  +7.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.19"
$ws.Range("E42").Value = "  +12.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.107"
$ws.Range("E44").Value = "  -3.56%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.06"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.24"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.40"
$ws.Range("E51").Value = "  +15.65%  "
